$d = $word.ActiveDocument

# The change trims the w:docDefaults block in styles.xml down to only the
# handful of properties that aren't already implied by Word's built-in
# defaults (rFonts/sz/szCs/lang for rPrDefault, spacing line/lineRule for
# pPrDefault). There's no direct Styles/Font/ParagraphFormat OM surface for
# docDefaults (writes there land on the individual style's own rPr/pPr), so
# we round-trip the whole package through WordOpenXML and rewrite the
# docDefaults fragment with a literal string replace.

$xml = $d.WordOpenXML

$oldDocDefaults = '<w:docDefaults><w:rPrDefault><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b w:val="0"/><w:i w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:color w:val="000000"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="none"/><w:shd w:val="clear" w:fill="auto"/><w:vertAlign w:val="baseline"/><w:lang w:val="en"/></w:rPr></w:rPrDefault><w:pPrDefault><w:pPr><w:keepNext w:val="0"/><w:keepLines w:val="0"/><w:widowControl/><w:pBdr><w:top w:val="nil" w:sz="0" w:space="0"/><w:left w:val="nil" w:sz="0" w:space="0"/><w:bottom w:val="nil" w:sz="0" w:space="0"/><w:right w:val="nil" w:sz="0" w:space="0"/><w:between w:val="nil" w:sz="0" w:space="0"/></w:pBdr><w:shd w:val="clear" w:fill="auto"/><w:spacing w:before="0" w:after="0" w:line="276" w:lineRule="auto"/><w:ind w:left="0" w:right="0" w:firstLine="0"/><w:contextualSpacing w:val="0"/><w:jc w:val="left"/></w:pPr></w:pPrDefault></w:docDefaults>'

$newDocDefaults = '<w:docDefaults><w:rPrDefault><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en"/></w:rPr></w:rPrDefault><w:pPrDefault><w:pPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr></w:pPrDefault></w:docDefaults>'

if ($xml.Contains($oldDocDefaults)) {
    $xml = $xml.Replace($oldDocDefaults, $newDocDefaults)
    $d.WordOpenXML = $xml
} else {
    Write-Output "docDefaults block not found verbatim; no change applied"
}
